# The author re-uploaded the deck after pruning three "placeholder" slides
# (empty title + generic workflow description in the content placeholder)
# that were left over from the outline template. The remaining slides keep
# their exact content/order - only the three entries with slide IDs
# 260, 261, 262 (originally at positions 8, 10 and 12) are removed from the
# slide list.
#
# Delete from the highest index down so earlier indices stay valid while we
# iterate.

$p = $ppt.ActivePresentation

$idsToRemove = @(260, 261, 262)

# Snapshot slide indices whose SlideID is in the removal set, highest first.
$indexesToDelete = @()
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($idsToRemove -contains $slide.SlideID) {
        $indexesToDelete += $i
    }
}
$indexesToDelete = $indexesToDelete | Sort-Object -Descending

foreach ($idx in $indexesToDelete) {
    $p.Slides.Item($idx).Delete()
}

Write-Output ("Slides remaining: " + $p.Slides.Count)
